# Apply crypto price/coin list update (GitHub Actions scraper run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.541'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05636'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.382'
$ws.Range("E6").Value = '5GateTokenGT'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.473'
$ws.Range("E7").Value = '6KuCoinTokenKCS'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8057'
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.058'
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1433'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07317'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03194'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02929'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09251'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001661'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.197'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04724'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006275'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0001503'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.968'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.137'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.01174'
$ws.Range("E24").Value = '23OneONE'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04167'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006921'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003506'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1039'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009843'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6811'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02239'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002103'
